$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-case row (row 5): CP-002 - Crear un turno valido
$ws.Range("B5").Value = "CP-002"
$ws.Range("C5").Value = 45595
$ws.Range("D5").Value = "Crear un turno valido"
$ws.Range("E5").Value = "Aprobado"
$ws.Range("F5").Value = "Funciono Correctamente"

# Match formatting of the first data row: centered "Caso de Prueba" column,
# and a real date format on the new Fecha cell.
$ws.Range("B5:B13").HorizontalAlignment = -4108
$ws.Range("C5").NumberFormat = "mm-dd-yy"

# Column widths grow to fit the new, longer content (target stored widths are
# 17.5546875 / 8.77734375 / 21.21875; ColumnWidth is specified net of the
# fixed ~0.8333 char padding Excel bakes into the stored <col width>).
$ws.Columns("D").ColumnWidth = 17.5546875 - 0.8333333333333334
$ws.Columns("E").ColumnWidth = 8.77734375 - 0.8333333333333334
$ws.Columns("F").ColumnWidth = 21.21875 - 0.8333333333333334
